# Deploying to gh-pages — add the "2022" column (S) to the stats table,
# mirroring the formatting already used by the "2021" column (R).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clone column R's formatting (number formats, fonts, borders, etc.)
#    into column S for the same row range, so new cells inherit the
#    correct visual style before we touch any values.
$ws.Range("R3:R33").Copy() | Out-Null
$ws.Range("S3:S33").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 2) Write the "2022" data into column S.
$values = @{
    3  = 2022
    4  = 1.9210869108320343
    5  = 1.020872301352429
    6  = 2.8415499553180767
    7  = 1.5924017665043597
    8  = 2.5011433798307796
    9  = 0.70098698968147144
    10 = 2.2312343573160249
    11 = 2.4764236727529938
    12 = 1.9888745417939038
    13 = 1.3057776932131271
    14 = 2.6056788910230639
    15 = 0
    16 = 0.65058422463372112
    17 = 0.65686622262510019
    18 = 0.64442124527961442
    19 = 2.5553368555544047
    20 = 1.807815324711445
    21 = 3.2928586128833093
    22 = 1.8387963974300983
    23 = 2.2260807622100529
    24 = 1.4582467499325562
    25 = 1.2245886088767601
    26 = 1.3105423773238725
    27 = 1.1375464261135158
    28 = 2.4791112740241377
    29 = 2.4279584268771761
    30 = 2.5408788313520994
    31 = 1.1238322680339958
    32 = 0.57553956834532372
    33 = 1.6467682173734046
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 19).Value = $values[$row]
}

# 3) The category-header rows (7, 10, 13, 16, 19, 22, 25, 28, 31) use a
#    bold variant of the regular numeric style in every other year
#    column — match that in the new column too.
foreach ($row in 7, 10, 13, 16, 19, 22, 25, 28, 31) {
    $ws.Cells.Item($row, 19).Font.Bold = $true
}

# 4) Move the active selection to T3, matching where the cursor landed
#    after the new column was filled in.
$ws.Range("T3").Select() | Out-Null
